# break fix from MBHS-1289
# The "Mapping" sheet's sample/test row data is being corrected:
#   - Row 2: Code changes from 1111111 -> 76543, and the Description
#     ("case type desc") is replaced with a real QA note.
#   - Row 3: the leftover test Code/Description (123 / "sssss") are
#     cleared out entirely, leaving only the Mapping value.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mapping")

$ws.Range("A2").Value = 76543
$ws.Range("B2").Value = "playwright, aqa test for upload mapping"

$ws.Range("A3").ClearContents()
$ws.Range("B3").ClearContents()

$ws.Activate() | Out-Null
$ws.Range("B3").Select() | Out-Null
